# Add a "Save" column (H) to the s_vals sheet.
# H1 gets the header "Save" (same bold/bordered style as the other headers),
# and H2:H56 get 0/1 flags derived from the existing "sum" column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - reuse the style already applied to G1 ("sum") so it matches
# the look of the rest of the header row (bold font, thin border box,
# centered/top-aligned).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 0
    39 = 1
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 0
    45 = 1
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 1
    54 = 1
    55 = 0
    56 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
